# The workbook's single sheet has a "Förändrad" (Changed) date column in C.
# Every data row (C2:C205) previously held the serial date 45181
# (2023-09-12) and the commit bumps it by one day to 45182 (2023-09-13).
# Walk the used range of column C and update any cell holding the old
# date to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 205
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ([Math]::Abs($current - 45181) -lt 0.0001) {
        $cell.Value = 45182
    }
}
